# Highlights "Pagina 3." in yellow (matching the other "Pagina N." cues in
# the document) and, in doing so, splits the run that used to hold
# " 3. Crear una agenda de negocio como la practica 5, en donde coloque la "
# into a highlighted " 3." run and a plain-formatted remainder run.

$d = $word.ActiveDocument

# Anchor on the unique sentence that starts this paragraph so we don't
# touch any of the other "Pagina N." occurrences in the document.
$anchor = $d.Content
$found = $anchor.Find.Execute("Pagina 3. Crear una agenda de negocio como la practica 5", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the target 'Pagina 3. ...' paragraph text"
}

$start = $anchor.Start

# "Pagina" is 6 characters.
$paginaRange = $d.Range($start, $start + 6)
$paginaRange.Font.HighlightColorIndex = 7

# " 3." is the next 3 characters (space, '3', '.').
$numberRange = $d.Range($start + 6, $start + 9)
$numberRange.Font.HighlightColorIndex = 7

Write-Output "Highlighted: [$($paginaRange.Text)][$($numberRange.Text)]"
